# Adds a new "Leading zeroes" test row to the Advanced Value Binder sample
# sheet, mirroring PHPExcel's AdvancedValueBinder unit test fixture update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label cell in column A.
$ws.Range("A6").Value = "Leading zeroes:"

# Force a text number format BEFORE assigning the value so the leading
# zeroes in "0001234" are preserved instead of being parsed as the
# number 1234.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "0001234"

# Widen column B a bit to comfortably fit the new values.
$ws.Columns.Item(2).ColumnWidth = 13.166666666666666

# Leave the new cell selected, matching the saved view state.
$ws.Range("B6").Select()
